$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose numeric (all-zero) period columns get cleared to blank text
# cells, matching the already-blank "B"/"D" columns on the same rows.
$rows = @(64, 79)
$cols = @(3,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33)

foreach ($r in $rows) {
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        # Assigning a bare "" clears/removes the cell entirely. Using a
        # leading apostrophe forces Excel to store an empty *text* value
        # (same representation already used by the neighboring empty
        # text cells), then the style is reset so the quote-prefix
        # formatting doesn't stick around.
        $cell.Value = "'"
        $cell.Style = "Normal"
    }
}
